# Regenerate the "K" column (column G) values in the save-data sheet.
# The original commit regenerated this column (K = strikeouts) using a
# fresh random draw per row, so here we simply overwrite the previously
# stored values with the newly "calculated" ones (s_vals), row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 3
    9  = 0
    10 = 0
    11 = 1
    12 = 1
    13 = 3
    14 = 1
    15 = 1
    16 = 2
    17 = 0
    18 = 0
    19 = 2
    20 = 1
    21 = 2
    22 = 3
    23 = 0
    24 = 2
    25 = 1
    26 = 1
    27 = 0
    28 = 2
    29 = 1
    30 = 1
    31 = 2
    32 = 3
    33 = 0
    34 = 2
    35 = 3
    36 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
